$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- First paragraph: text replacement ---
# The paragraph originally holds two runs: "**ID__AFFARS_pgi_5333_topic_6__ID**"
# followed by a run containing a single trailing space. Replace that whole
# stretch with the new ID placeholder text and no trailing space, which
# collapses the paragraph down to a single run.
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5333_topic_6__ID** ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5333__ID**", 2)

# --- First paragraph: formatting ---
# Indent left by 225 twips (=11.25 points).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (all four edges) with a 5-twip space/gap and no
# explicit line — matches <w:pBdr><w:top w:space="5"/> .../w:pBdr>.
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
